$d = $word.ActiveDocument

# 1. Remove the "_GoBack" bookmark from its original location (just before the
#    "This sample is set up to require Visual Studio 2017..." paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Update the "This sample is compatible with the ..." paragraph: merge the
#    two runs into a single run with the new SDK text, and restyle the
#    paragraph/run (Heading1 style + direct-formatting overrides), matching
#    what Word produces when the paragraph mark is re-typed along with the
#    following (now Heading1) paragraph mark.
$pCompat = $d.Paragraphs(2)
$rCompat = $pCompat.Range
$compatXml = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:pStyle w:val="Heading1"/>
    <w:spacing w:before="0"/>
    <w:rPr>
      <w:rFonts w:eastAsiaTheme="minorHAnsi" w:cs="Times New Roman"/>
      <w:i/>
      <w:color w:val="auto"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="22"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsiaTheme="minorHAnsi" w:cs="Times New Roman"/>
      <w:i/>
      <w:color w:val="auto"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t>This sample is compatible with the Windows 10 April 2018 Update SDK (17134)</w:t>
  </w:r>
</w:p>
"@
$rCompat.InsertXML($compatXml)

# 3. Turn the following empty paragraph into the new home of the "_GoBack"
#    bookmark.
$d = $word.ActiveDocument
$pEmpty = $d.Paragraphs(3)
$rEmpty = $pEmpty.Range
$bookmarkXml = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
"@
$rEmpty.InsertXML($bookmarkXml)
